# Update cryptocurrency price and volume(1h) figures on Sheet1
# Auto-generated from the authoritative cell-level diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.908.82'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '2.909.66'
$ws.Range("E3").Value = '  -0.06%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '588.77'
$ws.Range("E5").Value = '  +0.53%  '
$ws.Range("D6").Value = '144.65'
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = '6.89'
$ws.Range("E9").Value = '  +3.26%  '
$ws.Range("D10").Value = '0.141'
$ws.Range("E10").Value = '  -2.06%  '
$ws.Range("E11").Value = '  -2.04%  '
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").Value = '33.37'
$ws.Range("E13").Value = '  -0.04%  '
$ws.Range("E14").Value = '  +0.02%  '
$ws.Range("D15").Value = '3.388.80'
$ws.Range("D16").Value = '60.826.71'
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '6.67'
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("D18").Value = '2.906.71'
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("D19").Value = '431.84'
$ws.Range("E19").Value = '  +0.91%  '
$ws.Range("D20").Value = '13.33'
$ws.Range("E20").Value = '  -1.54%  '
$ws.Range("D21").Value = '0.675'
$ws.Range("E21").Value = '  -0.81%  '
$ws.Range("D22").Value = '7.09'
$ws.Range("D23").Value = '81.45'
$ws.Range("E23").Value = '  +0.93%  '
$ws.Range("D24").Value = '10.81'
$ws.Range("E24").Value = '  +1.58%  '
$ws.Range("D25").Value = '2.18'
$ws.Range("E25").Value = '  -2.28%  '
$ws.Range("D26").Value = '11.76'
$ws.Range("E26").Value = '  -1.75%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("E28").Value = '  +5.32%  '
$ws.Range("E29").Value = '  -0.86%  '
$ws.Range("D30").Value = '6.94'
$ws.Range("E30").Value = '  -3.94%  '
$ws.Range("D31").Value = '26.46'
$ws.Range("E31").Value = '  -0.39%  '
$ws.Range("E32").Value = '  +1.92%  '
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("D34").Value = '0.0₃0864'
$ws.Range("E34").Value = '  -1.12%  '
$ws.Range("E35").Value = '  -1.01%  '
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("D37").Value = '2.99'
$ws.Range("E37").Value = '  -0.49%  '
$ws.Range("D38").Value = '1.97'
$ws.Range("E38").Value = '  -1.27%  '
$ws.Range("E39").Value = '  -3.84%  '
$ws.Range("D40").Value = '8.54'
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("D41").Value = '40.91'
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = '0.281'
$ws.Range("E42").Value = '  -5.43%  '
$ws.Range("D43").Value = '376.39'
$ws.Range("E43").Value = '  -0.89%  '
$ws.Range("E44").Value = '  -1.50%  '
$ws.Range("D45").Value = '2.695.14'
$ws.Range("D46").Value = '133.65'
$ws.Range("E46").Value = '  +0.94%  '
$ws.Range("D48").Value = '23.66'
$ws.Range("E48").Value = '  -3.29%  '
$ws.Range("D49").Value = '0.106'
$ws.Range("E49").Value = '  -0.68%  '
$ws.Range("E50").Value = '  -2.63%  '
$ws.Range("E51").Value = '  -0.76%  '
